# Add a new "15. Feedback" test-case section to the API endpoint tracker.
# Matches the author's apparent entry order: endpoint + method typed first
# for the new row, then the module label, then the next row top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - Feedback info endpoint
$ws.Range("B19").Value2 = "/get-feedback-info"
$ws.Range("C19").Value2 = "Get - Provide Feed Back Keywords and Feed back Emoji Agaist a PIN"
$ws.Range("A19").Value2 = "15. Feedback"

# Row 20 - Save customer feedback endpoint
$ws.Range("A20").Value2 = "16.."
$ws.Range("B20").Value2 = "/save-customer-feedback"
$ws.Range("C20").Value2 = "POST - Store All Feedback from the user"

# Re-use the existing "section header" shaded style (as seen on A8/B8, A11/B11, ...)
# by copying formats from A8 onto the new Module/Endpoint header cells, so no new
# style entries are introduced.
$ws.Range("A8").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)

# Match the author's final selection in the saved file.
$ws.Range("G15").Select()
